# Fix a bug where looking up the OCMS id was not done right:
# the "OTL type ID" cell in row 7 (C7) held a bogus/malformed value
# ("OB0F4K3"); it should be a properly formatted OTL type id.
# While correcting the data, the cell is also given a distinct
# monospaced font (Monaco, 11pt) so object IDs are easier to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C7")
$cell.Value = "OB00454"
$cell.Font.Size = 11
$cell.Font.Name = "Monaco"

# Leave the selection where the editor left it after making the fix.
$ws.Range("C10").Select()
